$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.861.46'
$ws.Range('E2').Value = '  -1.55%  '
$ws.Range('D3').Value = '1.829.59'
$ws.Range('E3').Value = '  -1.07%  '
$ws.Range('D5').Value = '311.26'
$ws.Range('E5').Value = '  -0.87%  '
$ws.Range('D6').Value = '1.008'
$ws.Range('D7').Value = '0.4579'
$ws.Range('E7').Value = '  -0.46%  '
$ws.Range('D8').Value = '0.3672'
$ws.Range('E8').Value = '  -0.89%  '
$ws.Range('D9').Value = '0.07165'
$ws.Range('E9').Value = '  -1.60%  '
$ws.Range('D10').Value = '0.8762'
$ws.Range('E10').Value = '  -0.84%  '
$ws.Range('D11').Value = '0.07812'
$ws.Range('E11').Value = '  +0.22%  '
$ws.Range('D12').Value = '19.45'
$ws.Range('E12').Value = '  -2.24%  '
$ws.Range('D13').Value = '1.879.25'
$ws.Range('E13').Value = '  +1.71%  '
$ws.Range('D14').Value = '5.321'
$ws.Range('E14').Value = '  -0.86%  '
$ws.Range('D15').Value = '6.349'
$ws.Range('E15').Value = '  -2.93%  '
$ws.Range('D16').Value = '87.33'
$ws.Range('E16').Value = '  -4.46%  '
$ws.Range('D17').Value = '1.009'
$ws.Range('E17').Value = '  +0.67%  '
$ws.Range('D18').Value = '0.000008706'
$ws.Range('E18').Value = '  -2.71%  '
$ws.Range('D19').Value = '1.008'
$ws.Range('E19').Value = '  +0.65%  '
$ws.Range('D20').Value = '26.899.42'
$ws.Range('E20').Value = '  -1.50%  '
$ws.Range('E21').Value = '  -2.06%  '
$ws.Range('D22').Value = '4.980'
$ws.Range('E22').Value = '  -2.61%  '
$ws.Range('B23').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C23').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D23').Value = '2.081.11'
$ws.Range('E23').Value = '  +1.37%  '
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').Value = '10.45'
$ws.Range('E24').Value = '  -0.67%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = '1.991'
$ws.Range('E25').Value = '  +3.74%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '151.57'
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '18.15'
$ws.Range('E27').Value = '  -1.19%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').Value = '2.000'
$ws.Range('E28').Value = '  -2.67%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').Value = '113.64'
$ws.Range('E29').Value = '  -1.97%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = '4.906'
$ws.Range('E30').Value = '  -3.44%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').Value = '0.08774'
$ws.Range('E31').Value = '  -0.55%  '
$ws.Range('B32').Value = 'HuobiToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D32').Value = '3.106'
$ws.Range('E32').Value = '  -0.58%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = '0.7405'
$ws.Range('E33').Value = '  -3.97%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '4.472'
$ws.Range('E34').Value = '  -0.42%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '1.129'
$ws.Range('E35').Value = '  -3.28%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').Value = '2.482'
$ws.Range('E36').Value = '  -6.52%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').Value = '1.082'
$ws.Range('E37').Value = '  +0.25%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.01936'
$ws.Range('E38').Value = '  -0.92%  '
$ws.Range('D39').Value = '0.05125'
$ws.Range('E39').Value = '  -1.81%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = '2.912'
$ws.Range('E40').Value = '  -1.35%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = '6.928'
$ws.Range('E41').Value = '  -0.98%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = '0.4949'
$ws.Range('E42').Value = '  -3.42%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').Value = '0.1590'
$ws.Range('E43').Value = '  -2.51%  '
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').Value = '8.247'
$ws.Range('E44').Value = '  -1.76%  '
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').Value = '1.008'
$ws.Range('E45').Value = '  +0.73%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = '0.4648'
$ws.Range('E46').Value = '  -3.15%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '10.10'
$ws.Range('E47').Value = '  -1.54%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = '102.98'
$ws.Range('E48').Value = '  +0.33%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = '1.588'
$ws.Range('E49').Value = '  -3.62%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.06061'
$ws.Range('E50').Value = '  -2.47%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '64.55'
$ws.Range('E51').Value = '  -1.34%  '
